$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2250.75
$ws.Range("I40").Value = 2001
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2001
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -1826
$ws.Range("N40").Value = -3350

$ws.Range("H43").Value = 5567255.5
$ws.Range("J43").Value = 7938651
$ws.Range("L43").Value = 7938651
$ws.Range("N43").Value = -7938789

$ws.Range("H64").Value = 3854.9
$ws.Range("J64").Value = 3865
$ws.Range("L64").Value = 3865
$ws.Range("N64").Value = -4361

$ws.Range("H67").Value = 3854.9
$ws.Range("J67").Value = 3865
$ws.Range("L67").Value = 3865
$ws.Range("N67").Value = -5581

$ws.Range("H76").Value = 5449.8
$ws.Range("I76").Value = 6125
$ws.Range("J76").Value = 4999.6665
$ws.Range("K76").Value = 6125
$ws.Range("L76").Value = 4999.6665
$ws.Range("M76").Value = -5810
$ws.Range("N76").Value = -5629.6665

$ws.Range("H79").Value = 5449.8
$ws.Range("I79").Value = 6125
$ws.Range("J79").Value = 4999.6665
$ws.Range("K79").Value = 6125
$ws.Range("L79").Value = 4999.6665
$ws.Range("M79").Value = -5033
$ws.Range("N79").Value = -7183.6665

$ws.Range("H86").Value = 3163.4546
$ws.Range("I86").Value = 3699.8
$ws.Range("J86").Value = 2716.5
$ws.Range("K86").Value = 3699.8
$ws.Range("L86").Value = 2716.5
$ws.Range("M86").Value = -2576.8
$ws.Range("N86").Value = -4962.5

$ws.Range("H89").Value = 3163.4546
$ws.Range("I89").Value = 3699.8
$ws.Range("J89").Value = 2716.5
$ws.Range("K89").Value = 18499
$ws.Range("L89").Value = 13582.5
$ws.Range("M89").Value = -12883
$ws.Range("N89").Value = -24814.5

$ws.Range("H98").Value = 4924.7856
$ws.Range("I98").Value = 6275.1
$ws.Range("J98").Value = 1549
$ws.Range("K98").Value = 6275.1
$ws.Range("L98").Value = 1549
$ws.Range("M98").Value = -4777.1
$ws.Range("N98").Value = -4545

$ws.Range("H122").Value = 4924.7856
$ws.Range("I122").Value = 6275.1
$ws.Range("J122").Value = 1549
$ws.Range("K122").Value = 18825.3
$ws.Range("L122").Value = 4647
$ws.Range("M122").Value = -16375.3
$ws.Range("N122").Value = -9547

$ws.Range("H129").Value = 778.0244
$ws.Range("J129").Value = 860.94116
$ws.Range("L129").Value = 2582.82348
$ws.Range("N129").Value = -12582.82348

$ws.Range("H132").Value = 6066732.5
$ws.Range("I132").Value = 7096097.5
$ws.Range("J132").Value = 19213.25
$ws.Range("K132").Value = 21288292.5
$ws.Range("L132").Value = 57639.75
$ws.Range("M132").Value = -21285762.5
$ws.Range("N132").Value = -62699.75

$ws.Range("H135").Value = 29412658
$ws.Range("I135").Value = 473.125
$ws.Range("J135").Value = 100001900
$ws.Range("K135").Value = 4258.125
$ws.Range("L135").Value = 900017100
$ws.Range("M135").Value = -1723.125
$ws.Range("N135").Value = -900022170

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 16.2
$ws.Range("I26").Value = 28.5
$ws.Range("J26").Value = 8
$ws.Range("K26").Value = 28.5
$ws.Range("L26").Value = 8
$ws.Range("M26").Value = 301.5
$ws.Range("N26").Value = -668

$ws.Range("H29").Value = 4250
$ws.Range("J29").Value = 500
$ws.Range("L29").Value = 500
$ws.Range("N29").Value = -1116

$ws.Range("H32").Value = 6449.27
$ws.Range("I32").Value = 5056.384
$ws.Range("K32").Value = 5056.384
$ws.Range("M32").Value = -4769.384

$ws.Range("H33").Value = 33334666
$ws.Range("J33").Value = 2000
$ws.Range("L33").Value = 2000
$ws.Range("N33").Value = -2658

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H37").Value = 28000
$ws.Range("J37").Value = 28000
$ws.Range("L37").Value = 28000
$ws.Range("N37").Value = -28546

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H41").Value = 4925.6
$ws.Range("I41").Value = 4028.4443
$ws.Range("J41").Value = 13000
$ws.Range("K41").Value = 4028.4443
$ws.Range("L41").Value = 13000
$ws.Range("M41").Value = -3614.4443
$ws.Range("N41").Value = -13828

$ws.Range("H45").Value = 2621.9092
$ws.Range("I45").Value = 3523.5
$ws.Range("J45").Value = 1540
$ws.Range("K45").Value = 3523.5
$ws.Range("L45").Value = 1540
$ws.Range("M45").Value = -3146.5
$ws.Range("N45").Value = -2294

$ws.Range("H61").Value = 100001480
$ws.Range("I61").Value = 142858220
$ws.Range("J61").Value = 2400
$ws.Range("K61").Value = 142858220
$ws.Range("L61").Value = 2400
$ws.Range("M61").Value = -142858008
$ws.Range("N61").Value = -2824

$ws.Range("H132").Value = 3013.2
$ws.Range("I132").Value = 2573.5789
$ws.Range("J132").Value = 3535.25
$ws.Range("K132").Value = 7720.736699999999
$ws.Range("L132").Value = 10605.75
$ws.Range("M132").Value = -5190.736699999999
$ws.Range("N132").Value = -15665.75

$ws.Range("H136").Value = 100001480
$ws.Range("I136").Value = 142858220
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 428574660
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -428572110
$ws.Range("N136").Value = -12300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 260
$ws.Range("J31").Value = 260
$ws.Range("L31").Value = 260
$ws.Range("N31").Value = -764

$ws.Range("H36").Value = 866.3333
$ws.Range("I36").Value = 1149.5
$ws.Range("J36").Value = 300
$ws.Range("K36").Value = 1149.5
$ws.Range("L36").Value = 300
$ws.Range("M36").Value = -615.5
$ws.Range("N36").Value = -1368

$ws.Range("H37").Value = 2420
$ws.Range("I37").Value = 533.3333
$ws.Range("J37").Value = 5250
$ws.Range("K37").Value = 533.3333
$ws.Range("L37").Value = 5250
$ws.Range("M37").Value = -396.3333
$ws.Range("N37").Value = -5524

$ws.Range("H39").Value = 14999
$ws.Range("J39").Value = 14999
$ws.Range("L39").Value = 14999
$ws.Range("N39").Value = -15777

$ws.Range("H46").Value = 4982.5
$ws.Range("J46").Value = 4982.5
$ws.Range("L46").Value = 4982.5
$ws.Range("N46").Value = -5578.5

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()

$ws.Range("H134").Value = 1998.3334
$ws.Range("I134").Value = 1747.75
$ws.Range("K134").Value = 5243.25
$ws.Range("M134").Value = -2708.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1323.8
$ws.Range("I31").Value = 1289.8723
$ws.Range("K31").Value = 1289.8723
$ws.Range("M31").Value = -994.8723

$ws.Range("H34").Value = 1323.8
$ws.Range("I34").Value = 1289.8723
$ws.Range("K34").Value = 1289.8723
$ws.Range("M34").Value = -1087.8723

$ws.Range("H58").Value = 6140
$ws.Range("I58").Value = 1184.1818
$ws.Range("J58").Value = 10033.857
$ws.Range("K58").Value = 1184.1818
$ws.Range("L58").Value = 10033.857
$ws.Range("M58").Value = -981.1818000000001
$ws.Range("N58").Value = -10439.857

$ws.Range("H136").Value = 6140
$ws.Range("I136").Value = 1184.1818
$ws.Range("J136").Value = 10033.857
$ws.Range("K136").Value = 3552.5454
$ws.Range("L136").Value = 30101.571
$ws.Range("M136").Value = -1002.5454
$ws.Range("N136").Value = -35201.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 12666.833
$ws.Range("I5").Value = 1001
$ws.Range("K5").Value = 1001
$ws.Range("M5").Value = -889

$ws.Range("H126").Value = 2315.3845
$ws.Range("I126").Value = 1870
$ws.Range("J126").Value = 2697.1428
$ws.Range("K126").Value = 5610
$ws.Range("L126").Value = 8091.428400000001
$ws.Range("M126").Value = -3140
$ws.Range("N126").Value = -13031.4284

$ws.Range("H132").Value = 3028.8064
$ws.Range("I132").Value = 3338.3125
$ws.Range("J132").Value = 2698.6667
$ws.Range("K132").Value = 10014.9375
$ws.Range("L132").Value = 8096.000100000001
$ws.Range("M132").Value = -7484.9375
$ws.Range("N132").Value = -13156.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6707
$ws.Range("I40").Value = 2250.25
$ws.Range("J40").Value = 11163.75
$ws.Range("K40").Value = 2250.25
$ws.Range("L40").Value = 11163.75
$ws.Range("M40").Value = -2114.25
$ws.Range("N40").Value = -11435.75

$ws.Range("H122").Value = 11906395
$ws.Range("I122").Value = 17858542
$ws.Range("J122").Value = 2100.5715
$ws.Range("K122").Value = 53575626
$ws.Range("L122").Value = 6301.7145
$ws.Range("M122").Value = -53573176
$ws.Range("N122").Value = -11201.7145

$ws.Range("H136").Value = 2200.8572
$ws.Range("I136").Value = 2067.6667
$ws.Range("K136").Value = 6203.000100000001
$ws.Range("M136").Value = -3653.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 424.64706
$ws.Range("I113").Value = 302.27274
$ws.Range("J113").Value = 649
$ws.Range("K113").Value = 906.81822
$ws.Range("L113").Value = 1947
$ws.Range("M113").Value = 1263.18178
$ws.Range("N113").Value = -6287

$ws.Range("H136").Value = 1939.7
$ws.Range("I136").Value = 1851.75
$ws.Range("J136").Value = 1998.3334
$ws.Range("K136").Value = 5555.25
$ws.Range("L136").Value = 5995.0002
$ws.Range("M136").Value = -3005.25
$ws.Range("N136").Value = -11095.0002
